$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "和而泰"
$ws.Cells.Item(2, 2).Value = "和而泰"
$ws.Cells.Item(2, 3).Value = "和而泰"
$ws.Cells.Item(3, 1).Value = "绝味食品"
$ws.Cells.Item(3, 2).Value = "绝味食品"
$ws.Cells.Item(3, 3).Value = "立讯精密"
$ws.Cells.Item(4, 1).Value = "立讯精密"
$ws.Cells.Item(4, 2).Value = "复旦复华"
$ws.Cells.Item(4, 3).Value = "绝味食品"
$ws.Cells.Item(5, 1).Value = "首开股份"
$ws.Cells.Item(5, 2).Value = "立讯精密"
$ws.Cells.Item(5, 3).Value = "山子高科"
$ws.Cells.Item(6, 1).Value = "凯美特气"
$ws.Cells.Item(6, 2).Value = "创意信息"
$ws.Cells.Item(6, 3).Value = "上海建工"
$ws.Cells.Item(7, 1).Value = "山子高科"
$ws.Cells.Item(7, 2).Value = "赣锋锂业"
$ws.Cells.Item(7, 3).Value = "卧龙电驱"
$ws.Cells.Item(8, 1).Value = "赣锋锂业"
$ws.Cells.Item(8, 2).Value = "思科瑞"
$ws.Cells.Item(8, 3).Value = "首开股份"
$ws.Cells.Item(9, 1).Value = "上海建工"
$ws.Cells.Item(9, 2).Value = "天富能源"
$ws.Cells.Item(9, 3).Value = "欧菲光"
$ws.Cells.Item(10, 1).Value = "福龙马"
$ws.Cells.Item(10, 2).Value = "山子高科"
$ws.Cells.Item(10, 3).Value = "凯美特气"
$ws.Cells.Item(11, 1).Value = "卧龙电驱"
$ws.Cells.Item(11, 2).Value = "山河智能"
$ws.Cells.Item(11, 3).Value = "赣锋锂业"
$ws.Cells.Item(12, 1).Value = "山河智能"
$ws.Cells.Item(12, 2).Value = "上海建工"
$ws.Cells.Item(12, 3).Value = "天普股份"
$ws.Cells.Item(13, 1).Value = "三花智控"
$ws.Cells.Item(13, 2).Value = "卧龙电驱"
$ws.Cells.Item(13, 3).Value = "山河智能"
$ws.Cells.Item(14, 1).Value = "欧菲光"
$ws.Cells.Item(14, 2).Value = "凯美特气"
$ws.Cells.Item(14, 3).Value = "云南旅游"
$ws.Cells.Item(15, 1).Value = "复旦复华"
$ws.Cells.Item(15, 2).Value = "东华软件"
$ws.Cells.Item(15, 3).Value = "均胜电子"
$ws.Cells.Item(16, 1).Value = "创意信息"
$ws.Cells.Item(16, 2).Value = "天通股份"
$ws.Cells.Item(16, 3).Value = "三花智控"
$ws.Cells.Item(17, 1).Value = "金发科技"
$ws.Cells.Item(17, 2).Value = "三花智控"
$ws.Cells.Item(17, 3).Value = "工业富联"
$ws.Cells.Item(18, 1).Value = "思科瑞"
$ws.Cells.Item(18, 2).Value = "东方财富"
$ws.Cells.Item(18, 3).Value = "金发科技"
$ws.Cells.Item(19, 1).Value = "天普股份"
$ws.Cells.Item(19, 2).Value = "金发科技"
$ws.Cells.Item(19, 3).Value = "福龙马"
$ws.Cells.Item(20, 1).Value = "天富能源"
$ws.Cells.Item(20, 2).Value = "欧菲光"
$ws.Cells.Item(20, 3).Value = "华胜天成"
$ws.Cells.Item(21, 1).Value = "赛微电子"
$ws.Cells.Item(21, 2).Value = "福龙马"
$ws.Cells.Item(21, 3).Value = "利欧股份"
